$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the existing two country codes (row 3 and row 4).
$ws.Range("A3").Value = "AD"
$ws.Range("A4").Value = "AF"

# New country codes to append below the existing data, continuing the
# alternating row-style banding already applied to A3 (style odd) / A4
# (style even).
$codes = @("AG", "AI", "AL", "AO", "AQ", "AR", "AS", "AX", "BZ", "DZ", "US")

$row = 5
foreach ($code in $codes) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $code

    # Alternate formatting source row (A3 pattern for odd rows starting
    # at 5, A4 pattern for even rows), matching the existing banding.
    if (($row % 2) -eq 1) {
        $src = $ws.Range("A3")
    } else {
        $src = $ws.Range("A4")
    }
    $src.Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null

    $wholeRow = $ws.Rows.Item($row)
    $wholeRow.RowHeight = 18

    $row = $row + 1
}

$ws.Range("C9").Select() | Out-Null
